$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price + 1h volume change columns),
# plus a Coin/Link/Price/Volume row swap for Maker <-> Bittensor (rows 43/44).
# Price column (D) values are forced to Text so strings like "0.999" or
# "2.47" are not auto-coerced to numbers by Excel; Style is reset to Normal
# afterwards so no stray number-format style is left on the cell.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.845.50'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.58%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.928.97'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.21%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.86%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '549.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.05'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +8.77%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.512'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.87%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.930.68'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.128'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '4.76'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.445'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000220'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.80'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.33%  '
$ws.Range('E15').Value = '  +2.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.410.49'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.08%  '
$ws.Range('E17').Value = '  +7.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.926.43'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '57.762.51'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '416.56'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.35'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.694'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.70%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.45'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +7.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.98'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '79.44'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.67%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.47'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.02'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.41'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.26'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.97'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0971'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.24%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.67'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.938'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.25%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.07'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0₃0694'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +12.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '48.26'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.68'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.60'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +10.99%  '
$ws.Range('E41').Value = '  +1.76%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0345'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.55%  '
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '374.27'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.11%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.688.35'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.82%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '123.51'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.236'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.54%  '
$ws.Range('E48').Value = '  +2.40%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.95'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.90'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.99'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.30%  '
